# Update column G ("K") values for rows 2-24 to reflect the regenerated
# save_data (switch from Strike# to K, with recalculated std/mean and s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 4
    4  = 7
    5  = 6
    6  = 4
    7  = 2
    8  = 4
    9  = 1
    10 = 2
    11 = 2
    12 = 2
    13 = 2
    14 = 2
    15 = 1
    16 = 1
    17 = 3
    18 = 1
    19 = 1
    20 = 1
    21 = 3
    22 = 5
    23 = 3
    24 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
